$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.375.14'
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").Value = '3.764.79'
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '614.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.76'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.91%  '
$ws.Range("D7").Value = '3.766.81'
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.527'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.49%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.167'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.55'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.34%  '
$ws.Range("E12").Value = '  -1.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.18'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000254'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("D15").Value = '4.398.77'
$ws.Range("E15").Value = '  +0.27%  '
$ws.Range("D16").Value = '3.770.89'
$ws.Range("E16").Value = '  +0.22%  '
$ws.Range("D17").Value = '69.503.00'
$ws.Range("E17").Value = '  -0.51%  '
$ws.Range("E18").Value = '  -2.41%  '
$ws.Range("E19").Value = '  -1.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '500.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.35'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.88%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.723'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.57'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.84'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.96'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.96'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000135'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.40%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  +0.99%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.88%  '
$ws.Range("E32").Value = '  +3.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.63'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.07%  '
$ws.Range("E34").Value = '  -1.37%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("E36").Value = '  +1.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.14'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.351'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.29%  '
$ws.Range("E39").Value = '  +4.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '471.52'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.07'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +12.35%  '
$ws.Range("E42").Value = '  -4.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '49.77'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '45.25'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.41%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.61'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.12%  '
$ws.Range("D46").Value = '2.954.60'
$ws.Range("E46").Value = '  -3.51%  '
$ws.Range("E47").Value = '  -0.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '139.08'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.79%  '
